$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il16"
$ws.Range("C2").Value = "Kcnj10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.243623333333334
$ws.Range("H2").Value = 12.73087
$ws.Range("I2").Value = 0.2469246453968972
$ws.Range("J2").Value = 0.2469246453968973
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09324
$ws.Range("N2").Value = 0.27972
$ws.Range("O2").Value = 0.03378127458009811
$ws.Range("P2").Value = 0.03378127458009811
$ws.Range("Q2").Value = 0.3956754396000001
$ws.Range("R2").Value = 3.561078956400001
$ws.Range("S2").Value = 0.008341429246745944
$ws.Range("T2").Value = 0.008341429246745946

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il16"
$ws.Range("C3").Value = "Kcnj10"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.243623333333334
$ws.Range("H3").Value = 12.73087
$ws.Range("I3").Value = 0.2469246453968972
$ws.Range("J3").Value = 0.2469246453968973
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.293427
$ws.Range("N3").Value = 6.880281
$ws.Range("O3").Value = 0.8309189963150005
$ws.Range("P3").Value = 0.8309189963150007
$ws.Range("Q3").Value = 9.732440330496667
$ws.Range("R3").Value = 87.59196297447001
$ws.Range("S3").Value = 0.2051743785186273
$ws.Range("T3").Value = 0.2051743785186274

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il16"
$ws.Range("C4").Value = "Kcnj10"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.243623333333334
$ws.Range("H4").Value = 12.73087
$ws.Range("I4").Value = 0.2469246453968972
$ws.Range("J4").Value = 0.2469246453968973
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.373442
$ws.Range("N4").Value = 1.120326
$ws.Range("O4").Value = 0.1352997291049013
$ws.Range("P4").Value = 0.1352997291049013
$ws.Range("Q4").Value = 1.584747184846667
$ws.Range("R4").Value = 14.26272466362
$ws.Range("S4").Value = 0.03340883763152401
$ws.Range("T4").Value = 0.03340883763152402

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il16"
$ws.Range("C5").Value = "Kcnj10"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.865491666666667
$ws.Range("H5").Value = 17.596475
$ws.Range("I5").Value = 0.3412966552647515
$ws.Range("J5").Value = 0.3412966552647516
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.09324
$ws.Range("N5").Value = 0.27972
$ws.Range("O5").Value = 0.03378127458009811
$ws.Range("P5").Value = 0.03378127458009811
$ws.Range("Q5").Value = 0.5468984430000001
$ws.Range("R5").Value = 4.922085987000001
$ws.Range("S5").Value = 0.01152943602476766
$ws.Range("T5").Value = 0.01152943602476766

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il16"
$ws.Range("C6").Value = "Kcnj10"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.865491666666667
$ws.Range("H6").Value = 17.596475
$ws.Range("I6").Value = 0.3412966552647515
$ws.Range("J6").Value = 0.3412966552647516
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.293427
$ws.Range("N6").Value = 6.880281
$ws.Range("O6").Value = 0.8309189963150005
$ws.Range("P6").Value = 0.8309189963150007
$ws.Range("Q6").Value = 13.45207695660833
$ws.Range("R6").Value = 121.068692609475
$ws.Range("S6").Value = 0.283589874238254
$ws.Range("T6").Value = 0.2835898742382542

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il16"
$ws.Range("C7").Value = "Kcnj10"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.865491666666667
$ws.Range("H7").Value = 17.596475
$ws.Range("I7").Value = 0.3412966552647515
$ws.Range("J7").Value = 0.3412966552647516
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.373442
$ws.Range("N7").Value = 1.120326
$ws.Range("O7").Value = 0.1352997291049013
$ws.Range("P7").Value = 0.1352997291049013
$ws.Range("Q7").Value = 2.190420938983333
$ws.Range("R7").Value = 19.71378845085
$ws.Range("S7").Value = 0.04617734500172976
$ws.Range("T7").Value = 0.04617734500172978

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Il16"
$ws.Range("C8").Value = "Kcnj10"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.123111999999999
$ws.Range("H8").Value = 18.369336
$ws.Range("I8").Value = 0.3562868663317164
$ws.Range("J8").Value = 0.3562868663317164
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.09324
$ws.Range("N8").Value = 0.27972
$ws.Range("O8").Value = 0.03378127458009811
$ws.Range("P8").Value = 0.03378127458009811
$ws.Range("Q8").Value = 0.5709189628799999
$ws.Range("R8").Value = 5.138270665919999
$ws.Range("S8").Value = 0.01203582446083442
$ws.Range("T8").Value = 0.01203582446083442

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Il16"
$ws.Range("C9").Value = "Kcnj10"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.123111999999999
$ws.Range("H9").Value = 18.369336
$ws.Range("I9").Value = 0.3562868663317164
$ws.Range("J9").Value = 0.3562868663317164
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.293427
$ws.Range("N9").Value = 6.880281
$ws.Range("O9").Value = 0.8309189963150005
$ws.Range("P9").Value = 0.8309189963150007
$ws.Range("Q9").Value = 14.042910384824
$ws.Range("R9").Value = 126.386193463416
$ws.Range("S9").Value = 0.2960455253725665
$ws.Range("T9").Value = 0.2960455253725666

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Il16"
$ws.Range("C10").Value = "Kcnj10"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.123111999999999
$ws.Range("H10").Value = 18.369336
$ws.Range("I10").Value = 0.3562868663317164
$ws.Range("J10").Value = 0.3562868663317164
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.373442
$ws.Range("N10").Value = 1.120326
$ws.Range("O10").Value = 0.1352997291049013
$ws.Range("P10").Value = 0.1352997291049013
$ws.Range("Q10").Value = 2.286627191504
$ws.Range("R10").Value = 20.579644723536
$ws.Range("S10").Value = 0.0482055164983154
$ws.Range("T10").Value = 0.04820551649831541

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Il16"
$ws.Range("C11").Value = "Kcnj10"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.9536773333333334
$ws.Range("H11").Value = 2.861032
$ws.Range("I11").Value = 0.05549183300663471
$ws.Range("J11").Value = 0.05549183300663472
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.09324
$ws.Range("N11").Value = 0.27972
$ws.Range("O11").Value = 0.03378127458009811
$ws.Range("P11").Value = 0.03378127458009811
$ws.Range("Q11").Value = 0.08892087456
$ws.Range("R11").Value = 0.8002878710400001
$ws.Range("S11").Value = 0.001874584847750078
$ws.Range("T11").Value = 0.001874584847750079

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Il16"
$ws.Range("C12").Value = "Kcnj10"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.9536773333333334
$ws.Range("H12").Value = 2.861032
$ws.Range("I12").Value = 0.05549183300663471
$ws.Range("J12").Value = 0.05549183300663472
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.293427
$ws.Range("N12").Value = 6.880281
$ws.Range("O12").Value = 0.8309189963150005
$ws.Range("P12").Value = 0.8309189963150007
$ws.Range("Q12").Value = 2.187189345554667
$ws.Range("R12").Value = 19.684704109992
$ws.Range("S12").Value = 0.04610921818555253
$ws.Range("T12").Value = 0.04610921818555255

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Il16"
$ws.Range("C13").Value = "Kcnj10"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.9536773333333334
$ws.Range("H13").Value = 2.861032
$ws.Range("I13").Value = 0.05549183300663471
$ws.Range("J13").Value = 0.05549183300663472
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.373442
$ws.Range("N13").Value = 1.120326
$ws.Range("O13").Value = 0.1352997291049013
$ws.Range("P13").Value = 0.1352997291049013
$ws.Range("Q13").Value = 0.3561431707146667
$ws.Range("R13").Value = 3.205288536432
$ws.Range("S13").Value = 0.007508029973332097
$ws.Range("T13").Value = 0.007508029973332098

Write-Output "Applied Natmi update for Il16-Kcnj10 (Dr Hou advice)"